# "se realizo mas excepciones" - agrega 4 filas nuevas (43-46) a la planilla de carga
# con datos adicionales de ejemplo (mismas columnas A:M que las filas existentes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columnas: A=ID, B=Numero de Documento, C=Nombre, D=Edad, E=Genero, F=Pais,
#           G=Ciudad, H=Productos, I=Categorias, J=Costo, K=Dia de la Semana,
#           L=Mes, M=Anio
# Estas columnas contienen, en el archivo original, valores que *parecen*
# numericos (IDs, documentos, edad, costo, mes, anio) pero que deben guardarse
# como texto (igual que el resto de la planilla). Para que Excel no los
# reinterprete como numeros se usa el prefijo de comilla simple, que es la
# forma estandar de forzar texto sin tocar el NumberFormat de la celda.
$textCols = @("A", "B", "D", "J", "L", "M")

$newRows = @(
    @{
        A = "82560"; B = "775571623"; C = "Laura"; D = "72"; E = "MUJER"
        F = "COLOMBIA"; G = "Medellín"; H = "Aspiradora||Pan"
        I = "ELECTRODOMESTICO||ALIMENTO"; J = "985349.2435178054"
        K = "Jueves"; L = "04"; M = "2023"
    },
    @{
        A = "43963"; B = "998705680"; C = "Isabel"; D = "74"; E = "MUJER"
        F = "CHILE"; G = "Viña del Mar"; H = "Cepillo de dientes"
        I = "COSMETICO"; J = "375302.89530345425"
        K = "Martes"; L = "01"; M = "2021"
    },
    @{
        A = "79471"; B = "578927681"; C = "Paula"; D = "40"; E = "MUJER"
        F = "AUSTRALIA"; G = "Brisbane"; H = "Batidora||Horno"
        I = "ELECTRODOMESTICO||ELECTRODOMESTICO"; J = "650772.6144401672"
        K = "Miércoles"; L = "01"; M = "2021"
    },
    @{
        A = "91054"; B = "714860768"; C = "Laura"; D = "51"; E = "MUJER"
        F = "AUSTRALIA"; G = "Melbourne"; H = "Smartwatch"
        I = "TECNOLOGIA"; J = "719233.2007133622"
        K = "Viernes"; L = "01"; M = "2022"
    }
)

$startRow = 43
$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    foreach ($col in $cols) {
        $value = $rowData[$col]
        if ($textCols -contains $col) {
            $value = "'" + $value
        }
        $ws.Range($col + $rowNum).Value2 = $value
    }
}
